# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates per the authoritative xml diff,
# sheet by sheet, row by row (H..N columns of the Leve profit tables).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4830.7896
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3350
$ws.Range("H111").Value = 4124.5
$ws.Range("J111").Value = 3994.5
$ws.Range("L111").Value = 11983.5
$ws.Range("N111").Value = -18117.5
$ws.Range("H132").Value = 28460.494
$ws.Range("I132").Value = 31083.18
$ws.Range("J132").Value = 5730.5557
$ws.Range("K132").Value = 93249.54000000001
$ws.Range("L132").Value = 17191.6671
$ws.Range("M132").Value = -90719.54000000001
$ws.Range("N132").Value = -22251.6671
$ws.Range("H137").Value = 2177.2354
$ws.Range("J137").Value = 2702.3809
$ws.Range("L137").Value = 8107.1427
$ws.Range("N137").Value = -13207.1427
$ws.Range("H140").Value = 93788.55499999999
$ws.Range("J140").Value = 93788.55499999999
$ws.Range("L140").Value = 93788.55499999999
$ws.Range("N140").Value = -104148.555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 987.61536
$ws.Range("I2").Value = 868.5306399999999
$ws.Range("K2").Value = 868.5306399999999
$ws.Range("M2").Value = -755.5306399999999
$ws.Range("H32").Value = 4077.9
$ws.Range("I32").Value = 2047.2169
$ws.Range("J32").Value = 13992.412
$ws.Range("K32").Value = 2047.2169
$ws.Range("L32").Value = 13992.412
$ws.Range("M32").Value = -1760.2169
$ws.Range("N32").Value = -14566.412
$ws.Range("H61").Value = 10781.939
$ws.Range("I61").Value = 14309.077
$ws.Range("J61").Value = 8489.299999999999
$ws.Range("K61").Value = 14309.077
$ws.Range("L61").Value = 8489.299999999999
$ws.Range("M61").Value = -14097.077
$ws.Range("N61").Value = -8913.299999999999
$ws.Range("H74").Value = 1788271.5
$ws.Range("I74").Value = 2120737
$ws.Range("J74").Value = 5046.8184
$ws.Range("K74").Value = 2120737
$ws.Range("L74").Value = 5046.8184
$ws.Range("M74").Value = -2119863
$ws.Range("N74").Value = -6794.8184
$ws.Range("H76").Value = 14991.333
$ws.Range("J76").Value = 14991.333
$ws.Range("L76").Value = 14991.333
$ws.Range("N76").Value = -15667.333
$ws.Range("H77").Value = 1788271.5
$ws.Range("I77").Value = 2120737
$ws.Range("J77").Value = 5046.8184
$ws.Range("K77").Value = 10603685
$ws.Range("L77").Value = 25234.092
$ws.Range("M77").Value = -10599317
$ws.Range("N77").Value = -33970.092
$ws.Range("H79").Value = 14991.333
$ws.Range("J79").Value = 14991.333
$ws.Range("L79").Value = 14991.333
$ws.Range("N79").Value = -17331.333
$ws.Range("H116").Value = 987.61536
$ws.Range("I116").Value = 868.5306399999999
$ws.Range("K116").Value = 868.5306399999999
$ws.Range("M116").Value = 1425.46936
$ws.Range("H136").Value = 10781.939
$ws.Range("I136").Value = 14309.077
$ws.Range("J136").Value = 8489.299999999999
$ws.Range("K136").Value = 42927.231
$ws.Range("L136").Value = 25467.9
$ws.Range("M136").Value = -40377.231
$ws.Range("N136").Value = -30567.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 987.61536
$ws.Range("I3").Value = 868.5306399999999
$ws.Range("K3").Value = 868.5306399999999
$ws.Range("M3").Value = -754.5306399999999
$ws.Range("H88").Value = 36871.3
$ws.Range("J88").Value = 36871.3
$ws.Range("L88").Value = 36871.3
$ws.Range("N88").Value = -37683.3
$ws.Range("H91").Value = 36871.3
$ws.Range("J91").Value = 36871.3
$ws.Range("L91").Value = 36871.3
$ws.Range("N91").Value = -39679.3
$ws.Range("H134").Value = 509905.94
$ws.Range("I134").Value = 593707.0600000001
$ws.Range("J134").Value = 7099.143
$ws.Range("K134").Value = 1781121.18
$ws.Range("L134").Value = 21297.429
$ws.Range("M134").Value = -1778586.18
$ws.Range("N134").Value = -26367.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 65513
$ws.Range("H111").Value = 80702
$ws.Range("J111").Value = 80702
$ws.Range("L111").Value = 80702
$ws.Range("N111").Value = -88882
$ws.Range("H134").Value = 2812.0557
$ws.Range("I134").Value = 1920.7273
$ws.Range("J134").Value = 12616.667
$ws.Range("K134").Value = 5762.1819
$ws.Range("L134").Value = 37850.001
$ws.Range("M134").Value = -3227.1819
$ws.Range("N134").Value = -42920.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 104982.836
$ws.Range("J37").Value = 104982.836
$ws.Range("L37").Value = 314948.508
$ws.Range("N37").Value = -315172.508
$ws.Range("H113").Value = 1987.8572
$ws.Range("I113").Value = 1227.8572
$ws.Range("J113").Value = 2367.8572
$ws.Range("K113").Value = 3683.5716
$ws.Range("L113").Value = 7103.571599999999
$ws.Range("M113").Value = -1513.5716
$ws.Range("N113").Value = -11443.5716
$ws.Range("H122").Value = 766.0714
$ws.Range("I122").Value = 697.625
$ws.Range("J122").Value = 857.3333
$ws.Range("K122").Value = 6278.625
$ws.Range("L122").Value = 7715.9997
$ws.Range("M122").Value = -3828.625
$ws.Range("N122").Value = -12615.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2296.2708
$ws.Range("I102").Value = 1630.1818
$ws.Range("J102").Value = 3761.6667
$ws.Range("K102").Value = 1630.1818
$ws.Range("L102").Value = 3761.6667
$ws.Range("M102").Value = -8.181800000000067
$ws.Range("N102").Value = -7005.6667
$ws.Range("H132").Value = 1225.7576
$ws.Range("I132").Value = 1094.0625
$ws.Range("J132").Value = 1576.9445
$ws.Range("K132").Value = 3282.1875
$ws.Range("L132").Value = 4730.833500000001
$ws.Range("M132").Value = -752.1875
$ws.Range("N132").Value = -9790.833500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 69173.60000000001
$ws.Range("M74").ClearContents()
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("H77").Value = 69173.60000000001
$ws.Range("M77").ClearContents()
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("H92").Value = 60388.5
$ws.Range("J92").Value = 60388.5
$ws.Range("L92").Value = 60388.5
$ws.Range("N92").Value = -65380.5
$ws.Range("H122").Value = 4623.8213
$ws.Range("J122").Value = 5219.857
$ws.Range("L122").Value = 15659.571
$ws.Range("N122").Value = -20559.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 904.4
$ws.Range("I81").Value = 736
$ws.Range("K81").Value = 1472
$ws.Range("M81").Value = -411
$ws.Range("H84").Value = 904.4
$ws.Range("I84").Value = 736
$ws.Range("K84").Value = 7360
$ws.Range("M84").Value = -2056
$ws.Range("H122").Value = 2698.182
$ws.Range("I122").Value = 1967.4
$ws.Range("J122").Value = 4264.143
$ws.Range("K122").Value = 5902.200000000001
$ws.Range("L122").Value = 12792.429
$ws.Range("M122").Value = -3452.200000000001
$ws.Range("N122").Value = -17692.429
